$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreatureProto")

# New module_name values for column C ("creature module name"), set in the
# same order the original author entered them (rows 5,6,7,8,11,10,9,12,13,15,14)
# so the resulting shared-strings table is built in the same sequence.
$order = @(
    @{ Row = 5;  Value = "caocao" },
    @{ Row = 6;  Value = "caopi" },
    @{ Row = 7;  Value = "caoren" },
    @{ Row = 8;  Value = "daqiao" },
    @{ Row = 11; Value = "liubei" },
    @{ Row = 10; Value = "zhangfei" },
    @{ Row = 9;  Value = "zhugeliang" },
    @{ Row = 12; Value = "guojia" },
    @{ Row = 13; Value = "huanggai" },
    @{ Row = 15; Value = "jiaxu" },
    @{ Row = 14; Value = "jushou" }
)

foreach ($entry in $order) {
    $ws.Cells.Item($entry.Row, 3).Value = $entry.Value
}

# Apply the same cell format used by column B (rows 6-15) to C6:C15, matching
# the diff which adds style index 4 to those cells (C5 keeps its original,
# unstyled format).
$ws.Range("B6").Copy()
$ws.Range("C6:C15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to C15, as reflected in the diff.
$ws.Range("C15").Select()
